{"js": "// Add a new bullet-list item (\"Messages: Case Classes. State flows.\")\n// right after the \"Protocol: Context Driven Interaction REST P2P ...\" bullet.\n\nconst body = context.document.body;\n\n// Locate the paragraph that holds the \"Protocol: Context Driven Interaction\n// REST P2P ...\" bullet so the new paragraph lands in the right spot even if\n// the document layout shifts around a little.\nconst searchResults = body.search(\n  \"Protocol: Context Driven Interaction REST P2P\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nlet anchorParagraph;\nif (searchResults.items.length > 0) {\n  anchorParagraph = searchResults.items[0].paragraphs.getFirst();\n} else {\n  // Fallback: last non-empty paragraph in the body.\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n  const items = paragraphs.items;\n  let fallback = items[items.length - 1];\n  for (let i = items.length - 1; i >= 0; i--) {\n    if (items[i].text.trim().length > 0) {\n      fallback = items[i];\n      break;\n    }\n  }\n  anchorParagraph = fallback;\n}\n\n// Insert the new bullet paragraph right after the anchor; insertParagraph\n// clones the anchor paragraph's formatting (numbering, indentation,\n// borders/shading), which matches the existing list items.\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"Messages: Case Classes. State flows.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the bullet paragraph \"Protocol: Context Driven Interaction REST P2P ...\"\n# so the new bullet is inserted right after it, regardless of exact paragraph\n# index in the document.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Protocol: Context Driven Interaction REST P2P\")\n\nif ($found) {\n    # Expand the found text to the whole paragraph (wdParagraph = 4) so the\n    # range covers the paragraph mark too.\n    $rng.Expand(4) | Out-Null\n} else {\n    # Fallback: last non-empty paragraph in the document.\n    $paras = $d.Paragraphs\n    $anchor = $paras.Item($paras.Count)\n    for ($i = $paras.Count; $i -ge 1; $i--) {\n        $p = $paras.Item($i)\n        if ($p.Range.Text.Trim().Length -gt 0) {\n            $anchor = $p\n            break\n        }\n    }\n    $rng = $anchor.Range\n}\n\n$insertPos = $rng.End\n\n# Insert a new paragraph break after the anchor paragraph; this clones the\n# anchor's paragraph formatting (bullet numbering, indentation, borders,\n# shading), matching the rest of the list.\n$rng.InsertParagraphAfter()\n\n# Put the new bullet's text into the freshly created (still empty) paragraph.\n$newRange = $d.Range($insertPos, $insertPos)\n$newRange.InsertAfter(\"Messages: Case Classes. State flows.\")\n"}
